# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column value
$updates = @{
    2  = 153
    3  = 1706
    4  = 788
    5  = 1123
    6  = 33
    7  = 11934
    8  = 41
    11 = 407
    12 = 1112
    13 = 848
    14 = 13463
    15 = 13433
    20 = 282
    22 = 49
    23 = 95
    24 = 169
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
